# Reorders the attribute rows of the 3NF table (rows 63-73) on Sheet1.
# Each "before" row (a full A:G record, including its banding style) is
# relocated to a new row per the mapping below; row 67 ends up blank
# (the row element is removed entirely, matching the pre-existing gaps
# elsewhere in the sheet, e.g. row 62/39/31/...).
#
# Because several destinations overlap the original source rows, every
# source row is first snapshotted into a scratch area far below the
# used range, then written back out to its final destination, and the
# scratch area is cleared again at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# before-row -> after-row mapping for the A2:G73 "attribute" block
$moves = @{
    63 = 74
    64 = 68
    65 = 63
    66 = 71
    67 = 72
    68 = 73
    69 = 64
    70 = 70
    71 = 65
    72 = 66
    73 = 69
}

$scratchBase = 200

# 1) Snapshot every source row (A:G) into scratch rows 200-210, preserving
#    values and formatting (banding style travels with its row's data).
foreach ($src in $moves.Keys) {
    $scratchRow = $scratchBase + $src
    $ws.Range("A$src`:G$src").Copy($ws.Range("A$scratchRow`:G$scratchRow"))
}

# 2) Write each snapshot back out to its final destination row.
foreach ($src in $moves.Keys) {
    $dest = $moves[$src]
    $scratchRow = $scratchBase + $src
    $ws.Range("A$scratchRow`:G$scratchRow").Copy($ws.Range("A$dest`:G$dest"))
}

# 3) Row 67 is not a destination for anything -> it becomes the new gap.
$ws.Range("A67:G67").Clear()

# 4) Clean up the scratch area entirely.
foreach ($src in $moves.Keys) {
    $scratchRow = $scratchBase + $src
    $ws.Range("A$scratchRow`:G$scratchRow").Clear()
}

# 5) Restore the view state recorded in the saved workbook: scrolled so
#    row 34 is at the top, with G68 the active/selected cell.
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("G68").Select()
